$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Sheet1: "Review" section -------------------------------------------------
# Row 359 currently holds key "reviewTabHelp" / "Press the button and then share
# the link with your external reviewers:". It becomes the "get link" variant.
$ws.Range("B359").Value = "getReviewTabHelp"
$ws.Range("C359").Value = "Press the get link button and then share the link with your external reviewers:"

# Insert a brand-new row right below it for the "upload" variant, copying the
# formatting (style 45 on B:E) from the row we just edited.
$ws.Rows.Item(359).Copy()
$ws.Rows.Item(360).Insert()
$ws.Range("B360").Value = "updateReviewTabHelp"
$ws.Range("C360").Value = "Press the upload button and then share the link with your external reviewers:"

# --- Sheet1: new "loadingCourseError" row -------------------------------------
# This lands right after the "upgradeToStarterPlanToUseCommentsErrorMessage" row
# (originally row 366, now row 367 after the insertion above). Excel inserts a
# fresh unstyled row here (only B and C populated, no explicit style).
$ws.Rows.Item(368).Insert()
$ws.Range("B368").Value = "loadingCourseError"
$ws.Range("C368").Value = "Smth went wrong!"
